$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 186: Perejil, Primera, Vega Monumental Concepción, Región de Ñuble
$ws.Range("A186").Value = 11
$ws.Range("B186").Value = "Vega Monumental Concepción"
$ws.Range("C186").Value = "Bíobío"
$ws.Range("D186").Value = 44911
$ws.Range("D186").NumberFormat = $ws.Range("D185").NumberFormat
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 100112044
$ws.Range("G186").Value = "Perejil"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 200
$ws.Range("K186").Value = 700
$ws.Range("L186").Value = 800
$ws.Range("M186").Value = 750
$ws.Range("N186").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O186").Value = "Región de Ñuble"
$ws.Range("P186").Value = 750
$ws.Range("Q186").Value = 1
$ws.Range("R186").Value = "Hortaliza"

# New row 187: Perejil, Segunda, Vega Monumental Concepción, Región de Ñuble
$ws.Range("A187").Value = 11
$ws.Range("B187").Value = "Vega Monumental Concepción"
$ws.Range("C187").Value = "Bíobío"
$ws.Range("D187").Value = 44911
$ws.Range("D187").NumberFormat = $ws.Range("D185").NumberFormat
$ws.Range("E187").Value = 8
$ws.Range("F187").Value = 100112044
$ws.Range("G187").Value = "Perejil"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Segunda"
$ws.Range("J187").Value = 100
$ws.Range("K187").Value = 600
$ws.Range("L187").Value = 600
$ws.Range("M187").Value = 600
$ws.Range("N187").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O187").Value = "Región de Ñuble"
$ws.Range("P187").Value = 600
$ws.Range("Q187").Value = 1
$ws.Range("R187").Value = "Hortaliza"
